$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Auto2020-01-19-7076"
$ws.Range("A3").Value = "Auto2020-01-19-1149"
$ws.Range("A4").Value = "Auto2020-01-19-3940"
$ws.Range("A5").Value = "Auto2020-01-19-8243"
$ws.Range("A6").Value = "Auto2020-01-19-5764"
$ws.Range("A7").Value = "Auto2020-01-19-5911"
$ws.Range("A8").Value = "Auto2020-01-19-7859"
$ws.Range("A9").Value = "Auto2020-01-19-1605"
$ws.Range("A10").Value = "Auto2020-01-19-6059"
$ws.Range("A11").Value = "Auto2020-01-19-3893"
$ws.Range("A12").Value = "Auto2020-01-19-6856"
$ws.Range("A13").Value = "Auto2020-01-19-390"
$ws.Range("A14").Value = "Auto2020-01-19-3541"
$ws.Range("A15").Value = "Auto2020-01-19-4853"
$ws.Range("A16").Value = "Auto2020-01-19-5237"
$ws.Range("A17").Value = "Auto2020-01-19-6073"
$ws.Range("A18").Value = "Auto2020-01-19-1368"
$ws.Range("A19").Value = "Auto2020-01-19-8637"
$ws.Range("A20").Value = "Auto2020-01-19-6141"
$ws.Range("A21").Value = "Auto2020-01-19-9038"
$ws.Range("A22").Value = "Auto2020-01-19-7777"
$ws.Range("A23").Value = "Auto2020-01-19-6192"
$ws.Range("A24").Value = "Auto2020-01-19-6687"
$ws.Range("A25").Value = "Auto2020-01-19-132"
$ws.Range("A26").Value = "Auto2020-01-19-7413"
$ws.Range("A27").Value = "Auto2020-01-19-940"
$ws.Range("A28").Value = "Auto2020-01-19-8564"
$ws.Range("A29").Value = "Auto2020-01-19-2527"
$ws.Range("A30").Value = "Auto2020-01-19-7735"
$ws.Range("A31").Value = "Auto2020-01-19-9536"
$ws.Range("A32").Value = "Auto2020-01-19-4961"
$ws.Range("A33").Value = "Auto2020-01-19-2053"
$ws.Range("A34").Value = "Auto2020-01-19-7173"
$ws.Range("A35").Value = "Auto2020-01-19-682"
$ws.Range("A36").Value = "Auto2020-01-19-9241"
$ws.Range("A37").Value = "Auto2020-01-19-7378"
$ws.Range("A38").Value = "Auto2020-01-19-1759"
$ws.Range("A39").Value = "Auto2020-01-19-878"
$ws.Range("A40").Value = "Auto2020-01-19-2275"
$ws.Range("A41").Value = "Auto2020-01-19-8475"
$ws.Range("A42").Value = "Auto2020-01-19-2150"
$ws.Range("A43").Value = "Auto2020-01-19-947"
$ws.Range("A44").Value = "Auto2020-01-19-707"
$ws.Range("A45").Value = "Auto2020-01-19-1647"
$ws.Range("A46").Value = "Auto2020-01-19-1500"
$ws.Range("A47").Value = "Auto2020-01-19-1352"
$ws.Range("A48").Value = "Auto2020-01-19-6935"
$ws.Range("A49").Value = "Auto2020-01-19-3018"
$ws.Range("A50").Value = "Auto2020-01-19-9117"
$ws.Range("A51").Value = "Auto2020-01-19-5236"
$ws.Range("A52").Value = "Auto2020-01-19-5200"
$ws.Range("A53").Value = "Auto2020-01-19-627"
$ws.Range("A54").Value = "Auto2020-01-19-105"
$ws.Range("A55").Value = "Auto2020-01-19-2827"
$ws.Range("A56").Value = "Auto2020-01-19-5083"
$ws.Range("A57").Value = "Auto2020-01-19-8789"
$ws.Range("A58").Value = "Auto2020-01-19-2654"
$ws.Range("A59").Value = "Auto2020-01-19-7292"
$ws.Range("A60").Value = "Auto2020-01-19-5163"
$ws.Range("A61").Value = "Auto2020-01-19-7387"
$ws.Range("A62").Value = "Auto2020-01-19-9774"
$ws.Range("A63").Value = "Auto2020-01-19-9639"
$ws.Range("A64").Value = "Auto2020-01-19-6165"
$ws.Range("A65").Value = "Auto2020-01-19-964"
$ws.Range("A66").Value = "Auto2020-01-19-5650"
$ws.Range("A67").Value = "Auto2020-01-19-972"
$ws.Range("A68").Value = "Auto2020-01-19-3624"
$ws.Range("A69").Value = "Auto2020-01-19-2151"
$ws.Range("A70").Value = "Auto2020-01-19-1419"
$ws.Range("A71").Value = "Auto2020-01-19-6111"
$ws.Range("A72").Value = "Auto2020-01-19-7783"
$ws.Range("A73").Value = "Auto2020-01-19-5632"
$ws.Range("A74").Value = "Auto2020-01-19-9152"
$ws.Range("A75").Value = "Auto2020-01-19-759"
$ws.Range("A76").Value = "Auto2020-01-19-8404"
$ws.Range("A77").Value = "Auto2020-01-19-7415"
$ws.Range("A78").Value = "Auto2020-01-19-1881"
$ws.Range("A79").Value = "Auto2020-01-19-2102"
$ws.Range("A80").Value = "Auto2020-01-19-3738"
$ws.Range("A81").Value = "Auto2020-01-19-2386"
$ws.Range("A82").Value = "Auto2020-01-19-8951"
$ws.Range("A83").Value = "Auto2020-01-19-4914"
$ws.Range("A84").Value = "Auto2020-01-19-2695"
$ws.Range("A85").Value = "Auto2020-01-19-5553"
$ws.Range("A86").Value = "Auto2020-01-19-8253"
$ws.Range("A87").Value = "Auto2020-01-19-8450"
$ws.Range("A88").Value = "Auto2020-01-19-2276"
$ws.Range("A89").Value = "Auto2020-01-19-5975"
$ws.Range("A90").Value = "Auto2020-01-19-8463"
$ws.Range("A91").Value = "Auto2020-01-19-4149"
$ws.Range("A92").Value = "Auto2020-01-19-6883"
$ws.Range("A93").Value = "Auto2020-01-19-7416"
$ws.Range("A94").Value = "Auto2020-01-19-310"
$ws.Range("A95").Value = "Auto2020-01-19-8363"
$ws.Range("A96").Value = "Auto2020-01-19-2980"
$ws.Range("A97").Value = "Auto2020-01-19-1246"
$ws.Range("A98").Value = "Auto2020-01-19-1854"
$ws.Range("A99").Value = "Auto2020-01-19-3671"
$ws.Range("A100").Value = "Auto2020-01-19-8404"
$ws.Range("A101").Value = "Auto2020-01-19-4582"
$ws.Range("A102").Value = "Auto2020-01-19-7516"
$ws.Range("A103").Value = "Auto2020-01-19-4524"
$ws.Range("A104").Value = "Auto2020-01-19-8498"
$ws.Range("A105").Value = "Auto2020-01-19-9367"
$ws.Range("A106").Value = "Auto2020-01-19-7118"
$ws.Range("A107").Value = "Auto2020-01-19-6789"
$ws.Range("A108").Value = "Auto2020-01-19-3807"
$ws.Range("A109").Value = "Auto2020-01-19-5026"
$ws.Range("A110").Value = "Auto2020-01-19-561"
$ws.Range("A111").Value = "Auto2020-01-19-1457"
